# Update row 28 (2025Q2) metrics on Sheet1 of the recurrence metrics workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = 334
$ws.Range("D28").Value = 35
$ws.Range("E28").Value = 299
$ws.Range("F28").Value = 5.451713395638629
